$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.991.06'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.830.88'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.35'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6254'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.95%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.16'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07604'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2906'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.67'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.43%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07756'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.829.53'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.52%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.953'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6617'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.98%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.42'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009666'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +11.90%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.979'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.56%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.000.52'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '223.69'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.31'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.193'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.32'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1367'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.29%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.389'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.12%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.83'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.53%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.492'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.052'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.11%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.020'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.192'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05175'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.25%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.839'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.80%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7391'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.143'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.698'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.259.71'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.00%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.761'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01790'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.204'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8917'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.47'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.976.72'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000122'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '64.21'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.34%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5113'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3976'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.820'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05751'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.90%  '
